$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update index2Sequence value E7760 -> E7420 across K2:K37
$ws.Range("K2:K37").Value = "E7420"

# 2. Apply new font (Arial 11, black) to K2:K37 using a throwaway named style
#    so only a single new font/xf pair is minted.
$st = $wb.Styles.Add("TempFontStyle")
$st.Font.Size = 11
$st.Font.Color = 0
$ws.Range("K2:K37").Style = "TempFontStyle"
$wb.Styles.Item("TempFontStyle").Delete()

# 3. Convert L2:L37 boolean literals to formulas =FALSE()
for ($r = 2; $r -le 37; $r++) {
    $ws.Range("L$r").Formula = "=FALSE()"
}

# 4. Update sheet view / selection to match new scroll position + selection
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$ws.Range("K2:K37").Select() | Out-Null
